$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 13:42"

# Row 26
$ws.Range("B26").Value = 196989
$ws.Range("C26").Value = 2880
$ws.Range("D26").Value = 140652
$ws.Range("E26").Value = 48207
$ws.Range("G26").Value = 105
$ws.Range("H26").Value = 8130

# Row 38
$ws.Range("B38").Value = 90387
$ws.Range("C38").Value = 805
$ws.Range("D38").Value = 81037
$ws.Range("E38").Value = 8804
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = 546

# Row 57
$ws.Range("A57").Value = "Nepal"
$ws.Range("B57").Value = 47236
$ws.Range("C57").Value = 979
$ws.Range("D57").Value = 30677
$ws.Range("E57").Value = 16259
$ws.Range("G57").Value = 11
$ws.Range("H57").Value = 300

# Row 58
$ws.Range("A58").Value = "Costa Rica"
$ws.Range("B58").Value = 46920
$ws.Range("D58").Value = 18211
$ws.Range("E58").Value = 28231
$ws.Range("H58").Value = 478

# Row 59
$ws.Range("A59").Value = "Argelia"
$ws.Range("B59").Value = 46364
$ws.Range("D59").Value = 32745
$ws.Range("E59").Value = 12063
$ws.Range("H59").Value = 1556

# Row 62
$ws.Range("A62").Value = "Suiza"
$ws.Range("B62").Value = 44592
$ws.Range("C62").Value = 191
$ws.Range("D62").Value = 37100
$ws.Range("E62").Value = 5478
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 2014

# Row 63
$ws.Range("A63").Value = "Kirguistan"
$ws.Range("B63").Value = 44458
$ws.Range("C63").Value = 55
$ws.Range("D63").Value = 39960
$ws.Range("E63").Value = 3438
$ws.Range("H63").Value = 1060

# Row 73
$ws.Range("A73").Value = "Estado de Palestina"
$ws.Range("B73").Value = 26779
$ws.Range("C73").Value = 652
$ws.Range("D73").Value = 17270
$ws.Range("E73").Value = 9326
$ws.Range("G73").Value = 2
$ws.Range("H73").Value = 183

# Row 74
$ws.Range("A74").Value = "El Salvador"
$ws.Range("B74").Value = 26413
$ws.Range("C74").Value = 105
$ws.Range("D74").Value = 16137
$ws.Range("E74").Value = 9512
$ws.Range("G74").Value = 5
$ws.Range("H74").Value = 764

# Row 75
$ws.Range("A75").Value = "Australia"
$ws.Range("B75").Value = 26322
$ws.Range("C75").Value = 43
$ws.Range("D75").Value = 22603
$ws.Range("E75").Value = 2957
$ws.Range("G75").Value = 9
$ws.Range("H75").Value = 762

# Row 81
$ws.Range("A81").Value = "Libia"
$ws.Range("B81").Value = 18834
$ws.Range("C81").Value = 1085
$ws.Range("D81").Value = 2126
$ws.Range("E81").Value = 16412
$ws.Range("G81").Value = 11
$ws.Range("H81").Value = 296

# Row 82
$ws.Range("A82").Value = "Costa de Marfil"
$ws.Range("B82").Value = 18588
$ws.Range("D82").Value = 17472
$ws.Range("E82").Value = 997
$ws.Range("H82").Value = 119

# Row 83
$ws.Range("A83").Value = "Dinamarca"
$ws.Range("B83").Value = 17883
$ws.Range("D83").Value = 15760
$ws.Range("E83").Value = 1496
$ws.Range("H83").Value = 627

# Row 85
$ws.Range("B85").Value = 15352
$ws.Range("C85").Value = 33
$ws.Range("D85").Value = 14184
$ws.Range("E85").Value = 966
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 202

# Row 86
$ws.Range("B86").Value = 15127
$ws.Range("C86").Value = 37
$ws.Range("D86").Value = 12486
$ws.Range("E86").Value = 2018
$ws.Range("G86").Value = 6
$ws.Range("H86").Value = 623

# Row 146
$ws.Range("A146").Value = "Malta"
$ws.Range("B146").Value = 2076
$ws.Range("C146").Value = 37
$ws.Range("D146").Value = 1690
$ws.Range("H146").Value = 14

# Row 147
$ws.Range("A147").Value = "Sierra Leona"
$ws.Range("B147").Value = 2054
$ws.Range("D147").Value = 1611
$ws.Range("E147").Value = 372
$ws.Range("H147").Value = 71

# Row 165
$ws.Range("D165").Value = 853
$ws.Range("E165").Value = 161

# Row 218
$ws.Range("B218").Value = 9
$ws.Range("C218").Value = 3
$ws.Range("E218").Value = 4

